$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 9365.615
$ws.Range("I62").Value = 8940.429
$ws.Range("K62").Value = 8940.429
$ws.Range("M62").Value = -8316.429
# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 9365.615
$ws.Range("I65").Value = 8940.429
$ws.Range("K65").Value = 44702.145
$ws.Range("M65").Value = -41582.145
# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 36280.94
$ws.Range("I98").Value = 38585.332
$ws.Range("K98").Value = 38585.332
$ws.Range("M98").Value = -37087.332
# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 36280.94
$ws.Range("I122").Value = 38585.332
$ws.Range("K122").Value = 115755.996
$ws.Range("M122").Value = -113305.996
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 4346.54
$ws.Range("J138").Value = 5369.162
$ws.Range("L138").Value = 16107.486
$ws.Range("N138").Value = -26387.486
# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 9274.75
$ws.Range("I141").Value = 9663.362999999999
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 28990.089
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -23810.089
$ws.Range("N141").Value = -25360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 7125.25
$ws.Range("I5").Value = 15501
$ws.Range("K5").Value = 15501
$ws.Range("M5").Value = -15389
# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 2799.6
$ws.Range("I63").Value = 2249.5
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2249.5
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1563.5
$ws.Range("N63").Value = -6372
# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 2799.6
$ws.Range("I66").Value = 2249.5
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 11247.5
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -7815.5
$ws.Range("N66").Value = -31864
# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 7148192
$ws.Range("I97").Value = 6394.6665
$ws.Range("J97").Value = 28573584
$ws.Range("K97").Value = 6394.6665
$ws.Range("L97").Value = 28573584
$ws.Range("M97").Value = -5898.6665
$ws.Range("N97").Value = -28574576
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1855.3448
$ws.Range("I132").Value = 1427.6471
$ws.Range("J132").Value = 4971.4287
$ws.Range("K132").Value = 4282.9413
$ws.Range("L132").Value = 14914.2861
$ws.Range("M132").Value = -1752.9413
$ws.Range("N132").Value = -19974.2861

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 7125.25
$ws.Range("I4").Value = 15501
$ws.Range("K4").Value = 15501
$ws.Range("M4").Value = -15386
# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1287.25
$ws.Range("I107").Value = 1157.8889
$ws.Range("K107").Value = 1157.8889
$ws.Range("M107").Value = 762.1111000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 318.17648
$ws.Range("I7").Value = 244.83333
$ws.Range("J7").Value = 494.2
$ws.Range("K7").Value = 244.83333
$ws.Range("L7").Value = 494.2
$ws.Range("M7").Value = -131.83333
$ws.Range("N7").Value = -720.2
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 9733.333000000001
$ws.Range("I62").Value = 8377.615
$ws.Range("J62").Value = 11936.375
$ws.Range("K62").Value = 8377.615
$ws.Range("L62").Value = 11936.375
$ws.Range("M62").Value = -7753.615
$ws.Range("N62").Value = -13184.375
# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 9733.333000000001
$ws.Range("I65").Value = 8377.615
$ws.Range("J65").Value = 11936.375
$ws.Range("K65").Value = 41888.075
$ws.Range("L65").Value = 59681.875
$ws.Range("M65").Value = -38768.075
$ws.Range("N65").Value = -65921.875
# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 9471.857
$ws.Range("I107").Value = 11418.909
$ws.Range("K107").Value = 11418.909
$ws.Range("M107").Value = -9498.909
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 4370.4136
$ws.Range("I132").Value = 5045.5415
$ws.Range("K132").Value = 15136.6245
$ws.Range("M132").Value = -12606.6245
# Row 141: No Greater Treasure / Claro Walnut Necklace of Gathering
$ws.Range("H141").Value = 201894.56
$ws.Range("J141").Value = 216208.06
$ws.Range("L141").Value = 216208.06
$ws.Range("N141").Value = -226568.06

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80: Saucy for a Suitor / Hollandaise Sauce
$ws.Range("H80").Value = 62608384
$ws.Range("I80").Value = 250004500
$ws.Range("K80").Value = 750013500
$ws.Range("M80").Value = -750012564
# Row 83: Saved by the Sauce (L) / Hollandaise Sauce
$ws.Range("H83").Value = 62608384
$ws.Range("I83").Value = 250004500
$ws.Range("K83").Value = 2250040500
$ws.Range("M83").Value = -2250035820
# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 2094.9
$ws.Range("I107").Value = 612.5
$ws.Range("K107").Value = 1837.5
$ws.Range("M107").Value = 82.5
# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 4749.512
$ws.Range("J122").Value = 7670.9565
$ws.Range("L122").Value = 69038.6085
$ws.Range("N122").Value = -73938.6085
# Row 130: Blast from the Pasta / The Noodles of Elpis
$ws.Range("H130").Value = 16999.834
$ws.Range("I130").Value = 2000
$ws.Range("J130").Value = 19999.8
$ws.Range("K130").Value = 6000
$ws.Range("L130").Value = 59999.39999999999
$ws.Range("M130").Value = -980
$ws.Range("N130").Value = -70039.39999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57: Gold Is So Last Year / Electrum Circlet (Amber)
$ws.Range("H57").Value = 44999.5
$ws.Range("J57").Value = 44999.5
$ws.Range("L57").Value = 44999.5
$ws.Range("N57").Value = -46639.5
# Row 62: The Goggles, They Do Naught / Mythrite Goggles of Gathering
$ws.Range("H62").Value = 84999.5
$ws.Range("J62").Value = 84999.5
$ws.Range("L62").Value = 84999.5
$ws.Range("N62").Value = -86371.5
# Row 65: Peril Never Wore Safety Goggles (L) / Mythrite Goggles of Gathering
$ws.Range("H65").Value = 84999.5
$ws.Range("J65").Value = 84999.5
$ws.Range("L65").Value = 254998.5
$ws.Range("N65").Value = -261862.5
# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 39333
$ws.Range("I97").Value = 100000
$ws.Range("J97").Value = 8999.5
$ws.Range("K97").Value = 100000
$ws.Range("L97").Value = 8999.5
$ws.Range("M97").Value = -99504
$ws.Range("N97").Value = -9991.5
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2761.3489
$ws.Range("I132").Value = 2737.5405
$ws.Range("K132").Value = 8212.621500000001
$ws.Range("M132").Value = -5682.621500000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 47: Springtime for Coerthas / Boarskin Harness
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
# Row 51: Skirt Chaser / Boarskin Skirt
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 52: The Tao of Rabbits / Boarskin Harness
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
# Row 58: Handle with Care / Peisteskin Cesti
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
# Row 64: Glorified Hole-punchers / Archaeoskin Gloves of Aiming
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67: Treat Them with Kid Gloves (L) / Archaeoskin Gloves of Aiming
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 4707.16
$ws.Range("I68").Value = 2650
$ws.Range("J68").Value = 6323.5
$ws.Range("K68").Value = 2650
$ws.Range("L68").Value = 6323.5
$ws.Range("M68").Value = -1901
$ws.Range("N68").Value = -7821.5
# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 4707.16
$ws.Range("I71").Value = 2650
$ws.Range("J71").Value = 6323.5
$ws.Range("K71").Value = 13250
$ws.Range("L71").Value = 31617.5
$ws.Range("M71").Value = -9506
$ws.Range("N71").Value = -39105.5
# Row 76: Dragoon Drop Rate / Dhalmelskin Breeches of Maiming
$ws.Range("H76").Value = 9399
$ws.Range("I76").Value = 9399
$ws.Range("K76").Value = 9399
$ws.Range("M76").Value = -9061
# Row 79: Exploiting the Adroit (L) / Dhalmelskin Breeches of Maiming
$ws.Range("H79").Value = 9399
$ws.Range("I79").Value = 9399
$ws.Range("K79").Value = 9399
$ws.Range("M79").Value = -8229
# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 6539.7
$ws.Range("I100").Value = 2200
$ws.Range("J100").Value = 7624.625
$ws.Range("K100").Value = 2200
$ws.Range("L100").Value = 7624.625
$ws.Range("M100").Value = -1659
$ws.Range("N100").Value = -8706.625
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 7279.357
$ws.Range("I136").Value = 8999
$ws.Range("J136").Value = 6992.75
$ws.Range("K136").Value = 26997
$ws.Range("L136").Value = 20978.25
$ws.Range("M136").Value = -24447
$ws.Range("N136").Value = -26078.25
# Row 137: Lending Artisans a Hand / Br'aaxskin Halfgloves of Crafting
$ws.Range("H137").Value = 54199.332
$ws.Range("I137").Value = 39000
$ws.Range("J137").Value = 69398.664
$ws.Range("K137").Value = 39000
$ws.Range("L137").Value = 69398.664
$ws.Range("M137").Value = -33900
$ws.Range("N137").Value = -79598.664

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 273519.16
$ws.Range("I62").Value = 576166.5
$ws.Range("K62").Value = 576166.5
$ws.Range("M62").Value = -575542.5
# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 273519.16
$ws.Range("I65").Value = 576166.5
$ws.Range("K65").Value = 2880832.5
$ws.Range("M65").Value = -2877712.5
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 6364.379
$ws.Range("I132").Value = 7898.8296
$ws.Range("J132").Value = 2568.6316
$ws.Range("K132").Value = 23696.4888
$ws.Range("L132").Value = 7705.8948
$ws.Range("M132").Value = -21166.4888
$ws.Range("N132").Value = -12765.8948
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 336811.6
$ws.Range("I136").Value = 344162.97
$ws.Range("K136").Value = 1032488.91
$ws.Range("M136").Value = -1029938.91
